# Apply scraped market-price updates to each profession sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1248.3334
$ws.Range("I2").Value = 1248.3334
$ws.Range("K2").Value = 1248.3334
$ws.Range("M2").Value = -1135.3334
$ws.Range("H17").Value = 2250
$ws.Range("I17").Value = 2250
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 6750
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -6582
$ws.Range("N17").ClearContents()
$ws.Range("H38").Value = 29.625
$ws.Range("I38").Value = 29.625
$ws.Range("K38").Value = 88.875
$ws.Range("M38").Value = 283.125
$ws.Range("H40").Value = 4240
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H62").Value = 1979.6
$ws.Range("I62").Value = 1979.6
$ws.Range("K62").Value = 1979.6
$ws.Range("M62").Value = -1355.6
$ws.Range("H65").Value = 1979.6
$ws.Range("I65").Value = 1979.6
$ws.Range("K65").Value = 9898
$ws.Range("M65").Value = -6778
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H111").Value = 4798
$ws.Range("I111").Value = 6019.3335
$ws.Range("J111").Value = 2966
$ws.Range("K111").Value = 18058.0005
$ws.Range("L111").Value = 8898
$ws.Range("M111").Value = -14991.0005
$ws.Range("N111").Value = -15032
$ws.Range("H125").Value = 1517.5714
$ws.Range("I125").Value = 1256.3334
$ws.Range("J125").Value = 1987.8
$ws.Range("K125").Value = 11307.0006
$ws.Range("L125").Value = 17890.2
$ws.Range("M125").Value = -8847.000599999999
$ws.Range("N125").Value = -22810.2
$ws.Range("H138").Value = 7807.357
$ws.Range("J138").Value = 7217
$ws.Range("L138").Value = 21651
$ws.Range("N138").Value = -31931
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1986.6
$ws.Range("I45").Value = 1986.6
$ws.Range("K45").Value = 1986.6
$ws.Range("M45").Value = -1609.6
$ws.Range("H102").Value = 23335522
$ws.Range("I102").Value = 1431100.9
$ws.Range("J102").Value = 100001000
$ws.Range("K102").Value = 1431100.9
$ws.Range("L102").Value = 100001000
$ws.Range("M102").Value = -1429478.9
$ws.Range("N102").Value = -100004244
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 12666.333
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 12666.333
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 12666.333
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -13136.333
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H86").Value = 3608.5
$ws.Range("I86").Value = 3453.111
$ws.Range("K86").Value = 3453.111
$ws.Range("M86").Value = -2330.111
$ws.Range("H89").Value = 3608.5
$ws.Range("I89").Value = 3453.111
$ws.Range("K89").Value = 17265.555
$ws.Range("M89").Value = -11649.555
$ws.Range("H134").Value = 3267.75
$ws.Range("I134").Value = 3358
$ws.Range("J134").Value = 2997
$ws.Range("K134").Value = 10074
$ws.Range("L134").Value = 8991
$ws.Range("M134").Value = -7539
$ws.Range("N134").Value = -14061
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 717.3570999999999
$ws.Range("J22").Value = 600
$ws.Range("L22").Value = 600
$ws.Range("N22").Value = -1300
$ws.Range("H58").Value = 6400.2
$ws.Range("I58").Value = 1000
$ws.Range("J58").Value = 7750.25
$ws.Range("K58").Value = 1000
$ws.Range("L58").Value = 7750.25
$ws.Range("M58").Value = -797
$ws.Range("N58").Value = -8156.25
$ws.Range("H64").Value = 45000
$ws.Range("J64").Value = 45000
$ws.Range("L64").Value = 45000
$ws.Range("N64").Value = -45496
$ws.Range("H67").Value = 45000
$ws.Range("J67").Value = 45000
$ws.Range("L67").Value = 45000
$ws.Range("N67").Value = -46716
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H132").Value = 1470.3334
$ws.Range("I132").Value = 1206
$ws.Range("K132").Value = 3618
$ws.Range("M132").Value = -1088
$ws.Range("H134").Value = 1717.5
$ws.Range("I134").Value = 1623.3334
$ws.Range("K134").Value = 4870.0002
$ws.Range("M134").Value = -2335.0002
$ws.Range("H136").Value = 6400.2
$ws.Range("I136").Value = 1000
$ws.Range("J136").Value = 7750.25
$ws.Range("K136").Value = 3000
$ws.Range("L136").Value = 23250.75
$ws.Range("M136").Value = -450
$ws.Range("N136").Value = -28350.75
$ws.Range("H141").Value = 1036109.7
$ws.Range("J141").Value = 1036109.7
$ws.Range("L141").Value = 1036109.7
$ws.Range("N141").Value = -1046469.7
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 369.33334
$ws.Range("I50").Value = 369.33334
$ws.Range("K50").Value = 1108.00002
$ws.Range("M50").Value = -627.0000199999999
$ws.Range("H53").Value = 369.33334
$ws.Range("I53").Value = 369.33334
$ws.Range("K53").Value = 1108.00002
$ws.Range("M53").Value = -627.0000199999999
$ws.Range("H58").Value = 3440
$ws.Range("J58").Value = 3440
$ws.Range("L58").Value = 10320
$ws.Range("N58").Value = -10576
$ws.Range("H68").Value = 1097.75
$ws.Range("I68").Value = 3592
$ws.Range("J68").Value = 741.4286
$ws.Range("K68").Value = 10776
$ws.Range("L68").Value = 2224.2858
$ws.Range("M68").Value = -9965
$ws.Range("N68").Value = -3846.2858
$ws.Range("H71").Value = 1097.75
$ws.Range("I71").Value = 3592
$ws.Range("J71").Value = 741.4286
$ws.Range("K71").Value = 32328
$ws.Range("L71").Value = 6672.8574
$ws.Range("M71").Value = -28272
$ws.Range("N71").Value = -14784.8574
$ws.Range("H81").Value = 206.42857
$ws.Range("I81").Value = 206.42857
$ws.Range("K81").Value = 619.28571
$ws.Range("M81").Value = 503.71429
$ws.Range("H84").Value = 206.42857
$ws.Range("I84").Value = 206.42857
$ws.Range("K84").Value = 1857.85713
$ws.Range("M84").Value = 3758.14287
$ws.Range("H103").Value = 282.2
$ws.Range("J103").Value = 500
$ws.Range("L103").Value = 1500
$ws.Range("N103").Value = -3258
$ws.Range("H113").Value = 999.5
$ws.Range("I113").Value = 999
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 2997
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -827
$ws.Range("N113").Value = -7340
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 118.3
$ws.Range("I2").Value = 97.73333
$ws.Range("J2").Value = 180
$ws.Range("K2").Value = 97.73333
$ws.Range("L2").Value = 180
$ws.Range("M2").Value = 15.26667
$ws.Range("N2").Value = -406
$ws.Range("H58").Value = 47500
$ws.Range("I58").Value = 47500
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 47500
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -47223
$ws.Range("N58").ClearContents()
$ws.Range("H102").Value = 1924.0714
$ws.Range("I102").Value = 1378.7
$ws.Range("K102").Value = 1378.7
$ws.Range("M102").Value = 243.3
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 312.5
$ws.Range("I22").Value = 291.66666
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 291.66666
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = 3.333340000000021
$ws.Range("N22").Value = -1090
$ws.Range("H27").Value = 312.5
$ws.Range("I27").Value = 291.66666
$ws.Range("J27").Value = 500
$ws.Range("K27").Value = 291.66666
$ws.Range("L27").Value = 500
$ws.Range("M27").Value = -184.66666
$ws.Range("N27").Value = -714
$ws.Range("H40").Value = 1233.3334
$ws.Range("I40").Value = 1233.3334
$ws.Range("K40").Value = 1233.3334
$ws.Range("M40").Value = -1097.3334
$ws.Range("H68").Value = 4311.353
$ws.Range("I68").Value = 4589.3
$ws.Range("J68").Value = 3914.2856
$ws.Range("K68").Value = 4589.3
$ws.Range("L68").Value = 3914.2856
$ws.Range("M68").Value = -3840.3
$ws.Range("N68").Value = -5412.2856
$ws.Range("H71").Value = 4311.353
$ws.Range("I71").Value = 4589.3
$ws.Range("J71").Value = 3914.2856
$ws.Range("K71").Value = 22946.5
$ws.Range("L71").Value = 19571.428
$ws.Range("M71").Value = -19202.5
$ws.Range("N71").Value = -27059.428
$ws.Range("H105").Value = 29000
$ws.Range("J105").Value = 29000
$ws.Range("L105").Value = 29000
$ws.Range("N105").Value = -35988
$ws.Range("H122").Value = 3941.4
$ws.Range("I122").Value = 3902.6667
$ws.Range("K122").Value = 11708.0001
$ws.Range("M122").Value = -9258.000100000001
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7331.6665
$ws.Range("I62").Value = 6500
$ws.Range("K62").Value = 6500
$ws.Range("M62").Value = -5876
$ws.Range("H65").Value = 7331.6665
$ws.Range("I65").Value = 6500
$ws.Range("K65").Value = 32500
$ws.Range("M65").Value = -29380
$ws.Range("H136").Value = 3936.25
$ws.Range("I136").Value = 3875
$ws.Range("K136").Value = 11625
$ws.Range("M136").Value = -9075
